$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (D) and "Volumen" (J) values between row 2 and row 5.
$d2 = $ws.Range("D2").Value2
$d5 = $ws.Range("D5").Value2
$j2 = $ws.Range("J2").Value2
$j5 = $ws.Range("J5").Value2

$ws.Range("D2").Value = $d5
$ws.Range("D5").Value = $d2

$ws.Range("J2").Value = $j5
$ws.Range("J5").Value = $j2
